# Add a new book entry to the "Libros" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Libros")

# New row goes right after the last existing data row (row 18 -> row 19)
$newRow = 19

$ws.Cells.Item($newRow, 1).Value = 118
$ws.Cells.Item($newRow, 2).Value = "Mi Album de Viaje"
$ws.Cells.Item($newRow, 3).Value = "Javier"
$ws.Cells.Item($newRow, 4).Value = "Disponible"
